$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 (PED-B2-1 -> PED-B2-2, Session 2 -> 1, Date 31/08 -> 30/08) ---
$ws.Range("B3").Value = "PED-B2-2"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "30/08/2025"

# --- Clone the formatting of rows 2 and 3 onto the four new rows (4-7) ---
# Row 2's style pattern repeats on rows 4 and 6; row 3's style pattern repeats on rows 5 and 7.
$ws.Range("A2:G2").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Range("A6:G6").PasteSpecial(-4122)

$ws.Range("A3:G3").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$ws.Range("A7:G7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in values for the new rows ---
$ws.Range("A4").Value = "Year 5"
$ws.Range("B4").Value = "PED-B2-3"
$ws.Range("C4").Value = "pediatrics"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "30/08/2025"
$ws.Range("F4").Value = "10:00:00"
$ws.Range("G4").Value = 240

$ws.Range("A5").Value = "Year 5"
$ws.Range("B5").Value = "PED-B2-4"
$ws.Range("C5").Value = "pediatrics"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "30/08/2025"
$ws.Range("F5").Value = "10:00:00"
$ws.Range("G5").Value = 240

$ws.Range("A6").Value = "Year 5"
$ws.Range("B6").Value = "PED-B2-5"
$ws.Range("C6").Value = "pediatrics"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "30/08/2025"
$ws.Range("F6").Value = "10:00:00"
$ws.Range("G6").Value = 240

$ws.Range("A7").Value = "Year 5"
$ws.Range("B7").Value = "PED-B2-6"
$ws.Range("C7").Value = "pediatrics"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "30/08/2025"
$ws.Range("F7").Value = "10:00:00"
$ws.Range("G7").Value = 240
